# Edit: "Eetu Pihamäki" work-hours tracking workbook
# Adds a new logged work session (row 26) on the "Eetu Pihamäki" sheet,
# appends a trailing note/link to the existing row 24 task description,
# and updates the active selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# 1) Add the brand-new row 26 entry first (new shared string gets appended
#    before the modified row 24 string, matching the original authoring order).
$ws.Range("A26").Value = 43409
$ws.Range("B26").Value = 0.71875
$ws.Range("C26").Value = 0.81736111111111109
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = '1h 30 min Asentelin "Eclipse for Java and DSL Developers" pakettia ja yritin saada asennettua midPointin Log Viewer pluginia. Homma kaatui siihen, että Eclipse installerin help nappula ei toiminut, jonka kautta pluginin olisi voinut asentaa.  Myös pluginia itse kasattaessa tuli Maven erroria. Dokumentaatiossa kerrottiin myös, että SSL sertifikaatti ongelmia oli mahdollisesti tulossa. 30 min katselin midPointin dokumentaatiosta sekä Googlasin uusia connectroreita. Suunnittelin seuraavaan kertaan SSL/HTTPS -yhteyksien tekoa. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%205.11.2018.txt'

# 2) Append the follow-up note/link to the existing row 24 task text.
$ws.Range("F24").Value = '4 h midPointin ja testipalvelimen liittämistä niin että midPointista saadaan lisättyä käyttäjiä palvelimeen oikeilla asetuksilla (mm. admin käyttäjät saavat sudo oikeudet ja normaalikäyttäjät eivät pysty mm. muuttamaan asetustiedostoja jne.) Testattiin liittää useita erilaisia käyttäjiä palvelimeen midPointin kautta sekä ryhmä ja rooli jakoa käyttäjille. 1h 30 min Otin selvää midPointin lokeista /var/log/authlog, /var/log/syslog --> ei löytynyt mm. käyttäjien lisäämis tietoja ja ajankohtia. Katsoin myös midPoint GUI:sta lokeja jos löytyisi kiinnostavampaa loki tietoa, oli hieman epäselvää joten katsoin midPointin omasta dokumentaatiosta apua lokien selaamiseen --> https://wiki.evolveum.com/display/midPoint/Log+Viewer. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2031.10.2018.txt'

# 3) Resize the two rows whose wrapped text now needs more vertical space.
$ws.Rows.Item(24).RowHeight = 240
$ws.Rows.Item(26).RowHeight = 195

# 4) Update the sheet's active selection to reflect where the user ended up.
$ws.Activate()
$ws.Range("C26").Select()
